# Insert a new weekly record at row 81 (pushing existing rows 81..128 down to 82..129)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 81 - this shifts rows 81..128 down to 82..129
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new observation.
# Values mirror row 80 (same price point / metadata) but a newer date.
$ws.Cells.Item(81, 1).Value = 7
$ws.Cells.Item(81, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(81, 3).Value = "Ñuble"
$ws.Cells.Item(81, 4).Value = 44460
$ws.Cells.Item(81, 5).Value = 16
$ws.Cells.Item(81, 6).Value = 100112032
$ws.Cells.Item(81, 7).Value = "Zapallo italiano"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 120
$ws.Cells.Item(81, 11).Value = 11000
$ws.Cells.Item(81, 12).Value = 12000
$ws.Cells.Item(81, 13).Value = 11500
$ws.Cells.Item(81, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(81, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(81, 16).Value = 230
$ws.Cells.Item(81, 17).Value = 50
$ws.Cells.Item(81, 18).Value = "Hortaliza"
